$d = $word.ActiveDocument

# Locate the final (empty, en-GB) paragraph at the end of the document body.
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
$r = $lastPara.Range

# Build the OOXML package fragment with the four new paragraphs (title / body /
# title / body), each run + paragraph-mark carrying sz=28 (14pt) / szCs=28.
$xml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>Come guardare/commentare/valutare un video</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>L’utente raggiunge la piattaforma. Per guardare un video può scegliere se ricercare un video specifico tramite la barra di ricerca o selezionare uno dei video proposti all’arrivo sulla piattaforma di quelli già caricati da altri utenti. Una volta aperto il video ci saranno titolo e descrizione e i commenti lasciati dagli altri utenti registrati. In più se il video è piaciuto o meno c’è la possibilità di lasciare un like o un dislike e commentare il video stesso nella box apposita. Per commentare il video bisogna effettuare il login, se questo già è stato fatto allora ci sarà una box nel quale scrivere il commento ed inviarlo.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>Come segnalare un commento</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>Per segnalare un commento, esiste la funzione dal quale ogni utente registrato può raggiungere una volta aperto un video nell’apposita sezione “Commenti”. Da qui, come nel segnalare i video, si inserisce la motivazione della segnalazione che può essere una motivazione fra quelle di default o una scritta dall’utente stesso che vuole segnalare il commento. Una volta segnalato il commento si torna alla home della piattaforma.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# InsertXML on the (collapsed) final-paragraph range inserts the new
# paragraphs immediately before it, leaving that trailing paragraph in place.
$r.InsertXML($xml)

# The trailing paragraph (originally the document's only paragraph, carrying
# <w:lang w:val="en-GB"/>) must end up as a bare, empty paragraph. Select it
# and clear its character formatting so the leftover <w:rPr>/<w:lang> is gone.
$newCount = $d.Paragraphs.Count
$trailing = $d.Paragraphs.Item($newCount)
$trailing.Range.Select()
$word.Selection.ClearFormatting()

Write-Output "Done. ParagraphCount=$($d.Paragraphs.Count)"
